# Update the "courses" worksheet:
#  - department column (C) is renamed from the faculty name to the
#    simpler category labels used elsewhere ("Automotive" for the single
#    courses, "Packages" for the bundled package rows).
#  - promotionValidity column (R) values are cleared out (the promotion
#    text no longer applies), while keeping the existing cell formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2:C7 -> "Automotive", C8:C9 -> "Packages"
$ws.Range("C2:C7").Value = "Automotive"
$ws.Range("C8:C9").Value = "Packages"

# Clear the promotion validity text for every data row (R2:R9) but leave
# the cell style/formatting intact.
$ws.Range("R2:R9").ClearContents()
